$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.250.62"
$ws.Range("E2").Value = "  -5.72%  "
$ws.Range("D3").Value = "2.449.84"
$ws.Range("E3").Value = "  -8.46%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'534.55"
$ws.Range("E5").Value = "  -3.16%  "
$ws.Range("D6").Value = "'147.45"
$ws.Range("E6").Value = "  -6.71%  "
$ws.Range("E7").Value = "  -0.20%  "
$ws.Range("D8").Value = "'0.569"
$ws.Range("E8").Value = "  -3.93%  "
$ws.Range("D9").Value = "2.468.20"
$ws.Range("E9").Value = "  -7.81%  "
$ws.Range("E10").Value = "  -5.95%  "
$ws.Range("E11").Value = "  -2.47%  "
$ws.Range("D12").Value = "'5.34"
$ws.Range("E12").Value = "  +0.99%  "
$ws.Range("E13").Value = "  -4.54%  "
$ws.Range("D14").Value = "2.886.98"
$ws.Range("E14").Value = "  -8.24%  "
$ws.Range("E15").Value = "  -8.39%  "
$ws.Range("D16").Value = "59.229.98"
$ws.Range("E16").Value = "  -5.58%  "
$ws.Range("D17").Value = "'0.0000137"
$ws.Range("E17").Value = "  -6.17%  "
$ws.Range("D18").Value = "2.511.51"
$ws.Range("E18").Value = "  -6.04%  "
$ws.Range("D19").Value = "'11.11"
$ws.Range("E19").Value = "  -6.25%  "
$ws.Range("D20").Value = "'4.35"
$ws.Range("E20").Value = "  -5.57%  "
$ws.Range("D21").Value = "'323.28"
$ws.Range("E21").Value = "  -6.22%  "
$ws.Range("E22").Value = "  -3.15%  "
$ws.Range("D23").Value = "'5.73"
$ws.Range("E23").Value = "  -8.76%  "
$ws.Range("D24").Value = "'0.460"
$ws.Range("E24").Value = "  -8.87%  "
$ws.Range("D25").Value = "'60.44"
$ws.Range("E25").Value = "  -4.40%  "
$ws.Range("E26").Value = "  -4.34%  "
$ws.Range("E27").Value = "  -2.05%  "
$ws.Range("D28").Value = "'7.70"
$ws.Range("E28").Value = "  -5.71%  "
$ws.Range("E29").Value = "  -6.28%  "
$ws.Range("E30").Value = "  -6.31%  "
$ws.Range("D31").Value = "0.0₃0768"
$ws.Range("E31").Value = "  -10.49%  "
$ws.Range("D32").Value = "'1.24"
$ws.Range("E32").Value = "  -9.26%  "
$ws.Range("D34").Value = "'155.40"
$ws.Range("E34").Value = "  -6.41%  "
$ws.Range("E35").Value = "  -5.88%  "
$ws.Range("D36").Value = "'4.55"
$ws.Range("E36").Value = "  -6.34%  "
$ws.Range("E37").Value = "  -5.78%  "
$ws.Range("D38").Value = "'1.72"
$ws.Range("E38").Value = "  -3.21%  "
$ws.Range("D39").Value = "'5.86"
$ws.Range("E39").Value = "  -6.46%  "
$ws.Range("D40").Value = "'312.24"
$ws.Range("E40").Value = "  -10.30%  "
$ws.Range("D41").Value = "'36.72"
$ws.Range("E41").Value = "  -4.05%  "
$ws.Range("D42").Value = "'0.835"
$ws.Range("E42").Value = "  -12.70%  "
$ws.Range("E43").Value = "  -7.10%  "
$ws.Range("D44").Value = "'0.996"
$ws.Range("E44").Value = "  -0.26%  "
$ws.Range("E45").Value = "  -2.83%  "
$ws.Range("D46").Value = "'0.585"
$ws.Range("E46").Value = "  -5.00%  "
$ws.Range("E47").Value = "  -3.49%  "
$ws.Range("D48").Value = "'0.0525"
$ws.Range("E48").Value = "  -6.87%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'18.48"
$ws.Range("E49").Value = "  -9.06%  "
$ws.Range("B50").Value = "VeChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D50").Value = "'0.0228"
$ws.Range("E50").Value = "  -5.49%  "
$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D51").Value = "'18.82"
$ws.Range("E51").Value = "  -9.68%  "

$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
